# Apply updated "想去人数" (F) and "最低票价" (G) values to 展览 and 全部类型 sheets
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2084
$ws1.Range("G3").Value = 75
$ws1.Range("G4").Value = 60
$ws1.Range("F7").Value = 695
$ws1.Range("F10").Value = 42
$ws1.Range("F11").Value = 2542
$ws1.Range("F12").Value = 1614
$ws1.Range("F13").Value = 1590
$ws1.Range("F15").Value = 261
$ws1.Range("F16").Value = 638
$ws1.Range("F17").Value = 816
$ws1.Range("F19").Value = 325
$ws1.Range("F20").Value = 1094
$ws1.Range("F24").Value = 5457
$ws1.Range("F25").Value = 228
$ws1.Range("F26").Value = 798
$ws1.Range("F27").Value = 98
$ws1.Range("F31").Value = 228
$ws1.Range("F33").Value = 1055
$ws1.Range("F34").Value = 783
$ws1.Range("F39").Value = 1141
$ws1.Range("F41").Value = 109
$ws1.Range("F42").Value = 187
$ws1.Range("F44").Value = 83

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2084
$ws4.Range("G4").Value = 75
$ws4.Range("G5").Value = 60
$ws4.Range("F9").Value = 695
$ws4.Range("F14").Value = 42
$ws4.Range("F15").Value = 2542
$ws4.Range("F16").Value = 1614
$ws4.Range("F17").Value = 1590
$ws4.Range("F19").Value = 261
$ws4.Range("F20").Value = 638
$ws4.Range("F22").Value = 816
$ws4.Range("F24").Value = 325
$ws4.Range("F25").Value = 1094
$ws4.Range("F28").Value = 5457
$ws4.Range("F29").Value = 228
$ws4.Range("F30").Value = 798
$ws4.Range("F31").Value = 98
$ws4.Range("F35").Value = 228
$ws4.Range("F37").Value = 1055
$ws4.Range("F38").Value = 783
$ws4.Range("F41").Value = 1141
$ws4.Range("F43").Value = 109
$ws4.Range("F44").Value = 187
$ws4.Range("F46").Value = 83
